{"js": "// Update the date paragraph (first paragraph in the body).\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nif (paras.items.length > 0) {\n  paras.items[0].insertText(\"2023-07-02 Sunday\", \"Replace\");\n}\n\n// Update every multiplication answer cell in the single table, in row-major\n// order (20 rows x 5 columns = 100 cells), matching the diff's ordering.\nconst newValues = [\"36\u00d768=2448\", \"30\u00d732=960\", \"80\u00d763=5040\", \"19\u00d773=1387\", \"41\u00d713=533\", \"15\u00d741=615\", \"40\u00d727=1080\", \"85\u00d750=4250\", \"34\u00d753=1802\", \"50\u00d743=2150\", \"77\u00d775=5775\", \"58\u00d724=1392\", \"22\u00d710=220\", \"49\u00d757=2793\", \"10\u00d718=180\", \"29\u00d743=1247\", \"73\u00d711=803\", \"16\u00d752=832\", \"67\u00d724=1608\", \"80\u00d779=6320\", \"23\u00d732=736\", \"76\u00d745=3420\", \"55\u00d769=3795\", \"97\u00d727=2619\", \"64\u00d782=5248\", \"59\u00d784=4956\", \"72\u00d750=3600\", \"81\u00d775=6075\", \"16\u00d791=1456\", \"67\u00d732=2144\", \"10\u00d756=560\", \"66\u00d736=2376\", \"40\u00d7100=4000\", \"91\u00d740=3640\", \"29\u00d749=1421\", \"63\u00d788=5544\", \"86\u00d792=7912\", \"66\u00d769=4554\", \"77\u00d797=7469\", \"69\u00d715=1035\", \"29\u00d787=2523\", \"29\u00d719=551\", \"55\u00d711=605\", \"15\u00d766=990\", \"73\u00d712=876\", \"25\u00d719=475\", \"28\u00d776=2128\", \"49\u00d756=2744\", \"72\u00d720=1440\", \"11\u00d776=836\", \"61\u00d791=5551\", \"75\u00d776=5700\", \"81\u00d720=1620\", \"91\u00d754=4914\", \"60\u00d734=2040\", \"25\u00d736=900\", \"68\u00d746=3128\", \"46\u00d796=4416\", \"44\u00d761=2684\", \"94\u00d737=3478\", \"38\u00d728=1064\", \"82\u00d772=5904\", \"34\u00d750=1700\", \"60\u00d737=2220\", \"89\u00d743=3827\", \"70\u00d779=5530\", \"48\u00d768=3264\", \"69\u00d791=6279\", \"28\u00d766=1848\", \"16\u00d787=1392\", \"68\u00d780=5440\", \"89\u00d716=1424\", \"76\u00d793=7068\", \"22\u00d767=1474\", \"89\u00d799=8811\", \"21\u00d735=735\", \"43\u00d745=1935\", \"22\u00d778=1716\", \"39\u00d722=858\", \"22\u00d783=1826\", \"33\u00d748=1584\", \"84\u00d752=4368\", \"64\u00d734=2176\", \"53\u00d794=4982\", \"72\u00d739=2808\", \"37\u00d715=555\", \"81\u00d757=4617\", \"92\u00d734=3128\", \"29\u00d766=1914\", \"98\u00d734=3332\", \"39\u00d713=507\", \"20\u00d767=1340\", \"45\u00d793=4185\", \"38\u00d774=2812\", \"24\u00d798=2352\", \"69\u00d756=3864\", \"87\u00d710=870\", \"91\u00d787=7917\", \"46\u00d739=1794\", \"71\u00d773=5183\"];\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst rowCount = rows.items.length;\nconst colCount = newValues.length / rowCount;\n\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[i];\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date paragraph (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = '2023-07-02 Sunday'\n\n# Update every multiplication answer cell in the single table, in row-major\n# order (20 rows x 5 columns = 100 cells), matching the diff's ordering.\n$newValues = @(\n    '36\u00d768=2448',\n    '30\u00d732=960',\n    '80\u00d763=5040',\n    '19\u00d773=1387',\n    '41\u00d713=533',\n    '15\u00d741=615',\n    '40\u00d727=1080',\n    '85\u00d750=4250',\n    '34\u00d753=1802',\n    '50\u00d743=2150',\n    '77\u00d775=5775',\n    '58\u00d724=1392',\n    '22\u00d710=220',\n    '49\u00d757=2793',\n    '10\u00d718=180',\n    '29\u00d743=1247',\n    '73\u00d711=803',\n    '16\u00d752=832',\n    '67\u00d724=1608',\n    '80\u00d779=6320',\n    '23\u00d732=736',\n    '76\u00d745=3420',\n    '55\u00d769=3795',\n    '97\u00d727=2619',\n    '64\u00d782=5248',\n    '59\u00d784=4956',\n    '72\u00d750=3600',\n    '81\u00d775=6075',\n    '16\u00d791=1456',\n    '67\u00d732=2144',\n    '10\u00d756=560',\n    '66\u00d736=2376',\n    '40\u00d7100=4000',\n    '91\u00d740=3640',\n    '29\u00d749=1421',\n    '63\u00d788=5544',\n    '86\u00d792=7912',\n    '66\u00d769=4554',\n    '77\u00d797=7469',\n    '69\u00d715=1035',\n    '29\u00d787=2523',\n    '29\u00d719=551',\n    '55\u00d711=605',\n    '15\u00d766=990',\n    '73\u00d712=876',\n    '25\u00d719=475',\n    '28\u00d776=2128',\n    '49\u00d756=2744',\n    '72\u00d720=1440',\n    '11\u00d776=836',\n    '61\u00d791=5551',\n    '75\u00d776=5700',\n    '81\u00d720=1620',\n    '91\u00d754=4914',\n    '60\u00d734=2040',\n    '25\u00d736=900',\n    '68\u00d746=3128',\n    '46\u00d796=4416',\n    '44\u00d761=2684',\n    '94\u00d737=3478',\n    '38\u00d728=1064',\n    '82\u00d772=5904',\n    '34\u00d750=1700',\n    '60\u00d737=2220',\n    '89\u00d743=3827',\n    '70\u00d779=5530',\n    '48\u00d768=3264',\n    '69\u00d791=6279',\n    '28\u00d766=1848',\n    '16\u00d787=1392',\n    '68\u00d780=5440',\n    '89\u00d716=1424',\n    '76\u00d793=7068',\n    '22\u00d767=1474',\n    '89\u00d799=8811',\n    '21\u00d735=735',\n    '43\u00d745=1935',\n    '22\u00d778=1716',\n    '39\u00d722=858',\n    '22\u00d783=1826',\n    '33\u00d748=1584',\n    '84\u00d752=4368',\n    '64\u00d734=2176',\n    '53\u00d794=4982',\n    '72\u00d739=2808',\n    '37\u00d715=555',\n    '81\u00d757=4617',\n    '92\u00d734=3128',\n    '29\u00d766=1914',\n    '98\u00d734=3332',\n    '39\u00d713=507',\n    '20\u00d767=1340',\n    '45\u00d793=4185',\n    '38\u00d774=2812',\n    '24\u00d798=2352',\n    '69\u00d756=3864',\n    '87\u00d710=870',\n    '91\u00d787=7917',\n    '46\u00d739=1794',\n    '71\u00d773=5183'\n)\n\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n"}
